$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'322.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-2.44%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'15"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'42.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-5.98%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'15"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.213"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-4.72%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'15"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.08221"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.93%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'15"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'4.336"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-2.65%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'15"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'1.768"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-14.43%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'15"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.9498"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.50%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'15"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.1123"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.98%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'15"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.1890"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.95%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'15"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.09424"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-3.01%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'15"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.04638"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.87%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'15"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'7.452"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-20.92%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'15"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.1058"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.08%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'15"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.001303"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.38%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'15"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.005701"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-4.91%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'15"
$ws.Range("G16").Style = "Normal"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D17").Value = "'0.004280"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-3.78%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'15"
$ws.Range("G17").Style = "Normal"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.354"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.01%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'15"
$ws.Range("G18").Style = "Normal"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.563"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.53%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'15"
$ws.Range("G19").Style = "Normal"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3365"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.32%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'15"
$ws.Range("G20").Style = "Normal"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1391"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.47%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'15"
$ws.Range("G21").Style = "Normal"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2554"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.15%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'15"
$ws.Range("G22").Style = "Normal"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04169"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.13%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'15"
$ws.Range("G23").Style = "Normal"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001252"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-3.98%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'15"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001221"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-6.44%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'15"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.0002985"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.07%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'15"
$ws.Range("G26").Style = "Normal"
$ws.Range("G27").Value = "'15"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'15"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'15"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'15"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'15"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'15"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'15"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'15"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'15"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'15"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'15"
$ws.Range("G37").Style = "Normal"
$ws.Range("D38").Value = "'0.02646"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-3.77%"
$ws.Range("E38").Style = "Normal"
$ws.Range("G38").Value = "'15"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.05667"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.28%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'15"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.008147"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.65%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'15"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.1404"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.70%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'15"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.006567"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-9.72%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'15"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.001994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-6.12%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'15"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.007704"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.88%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'15"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.3496"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.25%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'15"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006772"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.45%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'15"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.21%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'15"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.003321"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-4.66%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'15"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.004109"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'15.83%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'15"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.21%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'15"
$ws.Range("G50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.21%"
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").Value = "'15"
$ws.Range("G51").Style = "Normal"
